$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -10.9133
$ws.Range("D3").Value = -7.247599999999992
$ws.Range("E8").Value = 16.5814
$ws.Range("E11").Value = 16.95689999999999
$ws.Range("B12").Value = 5.206599999999997
$ws.Range("C14").Value = -13.86759999999999
$ws.Range("E14").Value = 16.51590000000001
$ws.Range("E15").Value = 16.09270000000001
$ws.Range("C26").Value = -12.07280000000001
$ws.Range("D30").Value = -7.430400000000002
$ws.Range("C31").Value = -13.6905
$ws.Range("B32").Value = 6.306599999999999
$ws.Range("C35").Value = -12.91870000000001
$ws.Range("B36").Value = 8.853500000000006
$ws.Range("E36").Value = 16.14840000000001
$ws.Range("C37").Value = -13.4062
$ws.Range("B38").Value = 5.3337
$ws.Range("D44").Value = -7.230700000000001
$ws.Range("C45").Value = -13.45189999999999
$ws.Range("B46").Value = 6.490300000000001
$ws.Range("B54").Value = 4.618800000000002
$ws.Range("B55").Value = 5.339999999999996
$ws.Range("C57").Value = -14.44469999999999
$ws.Range("D58").Value = -8.113299999999995
$ws.Range("E64").Value = 17.5285
$ws.Range("B67").Value = 5.171799999999998
$ws.Range("B69").Value = 5.333399999999995
$ws.Range("B72").Value = 5.3447
$ws.Range("D84").Value = -8.360499999999998
$ws.Range("D89").Value = -7.202299999999995
$ws.Range("E89").Value = 17.77980000000002
$ws.Range("B91").Value = 5.523300000000002
$ws.Range("D91").Value = -6.245200000000001
$ws.Range("D92").Value = -6.085600000000002
$ws.Range("B99").Value = 4.643199999999998
$ws.Range("C100").Value = -12.8185
$ws.Range("C102").Value = -13.3348
$ws.Range("D102").Value = -7.948099999999998
